# Shorten the "Game Rules" slide (slide 2):
#  - "Sieger ist Spieler mit meisten Punkten am Ende"
#       -> "Sieger ist Spieler mit den meisten Punkten am Ende"
#  - "Spiel besteht aus 10 Runden à 20 Sekunden"
#       -> "Spiel besteht aus 10 Runden à 20 Sekunden plus einem Event."
#  - remove the two following bullets ("Am Ende jeder Runde ..." and
#    "Startpositionen nicht von Events betroffen")

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# --- Paragraph 3: add "den" ---------------------------------------------
# Setting .Text directly would keep the shared prefix/suffix as their own
# runs (PowerPoint's usual "reuse formatting" diffing). Routing the update
# through a disjoint placeholder first means the final assignment has no
# overlap with what's already there, so it lands back as a single run
# with the original rPr intact.
$para3 = $tr.Paragraphs(3, 1)
$para3.Text = "________________________________________________________"
$para3 = $tr.Paragraphs(3, 1)
$para3.Text = "Sieger ist Spieler mit den meisten Punkten am Ende"

# --- Paragraph 5: append "plus einem Event." ----------------------------
$para5 = $tr.Paragraphs(5, 1)
$para5.Text = "________________________________________________________"
$para5 = $tr.Paragraphs(5, 1)
$para5.Text = "Spiel besteht aus 10 Runden à 20 Sekunden plus einem Event."

# --- Remove the next two bullet paragraphs ------------------------------
# After each delete, the following paragraph slides into slot 6.
$para6 = $tr.Paragraphs(6, 1)
$para6.Delete()
$para6 = $tr.Paragraphs(6, 1)
$para6.Delete()
